$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '70.958.06'
$ws.Range('E2').Value = '  +2.93%  '

# Row 3
$ws.Range('D3').Value = '3.804.35'
$ws.Range('E3').Value = '  +1.09%  '

# Row 4
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').Value = '699.40'
$ws.Range('E5').Value = '  +8.29%  '

# Row 6
$ws.Range('D6').Value = '172.61'
$ws.Range('E6').Value = '  +4.39%  '

# Row 7
$ws.Range('D7').Value = '3.802.38'
$ws.Range('E7').Value = '  +1.16%  '

# Row 8
$ws.Range('E8').Value = '  -0.08%  '

# Row 9
$ws.Range('E9').Value = '  +1.30%  '

# Row 10
$ws.Range('E10').Value = '  +2.96%  '

# Row 11
$ws.Range('D11').Value = '7.38'
$ws.Range('E11').Value = '  +7.07%  '

# Row 12
$ws.Range('E12').Value = '  +1.15%  '

# Row 13
$ws.Range('E13').Value = '  +8.11%  '

# Row 14
$ws.Range('D14').Value = '36.44'
$ws.Range('E14').Value = '  +4.78%  '

# Row 15
$ws.Range('D15').Value = '4.446.19'
$ws.Range('E15').Value = '  +1.07%  '

# Row 16
$ws.Range('D16').Value = '3.821.69'
$ws.Range('E16').Value = '  +1.61%  '

# Row 17
$ws.Range('D17').Value = '70.934.40'
$ws.Range('E17').Value = '  +2.91%  '

# Row 18
$ws.Range('D18').Value = '17.88'
$ws.Range('E18').Value = '  +1.43%  '

# Row 19
$ws.Range('E19').Value = '  +3.17%  '

# Row 20
$ws.Range('E20').Value = '  +0.29%  '

# Row 21
$ws.Range('D21').Value = '11.10'
$ws.Range('E21').Value = '  +16.34%  '

# Row 22
$ws.Range('D22').Value = '482.08'
$ws.Range('E22').Value = '  +3.00%  '

# Row 23
$ws.Range('E23').Value = '  +1.63%  '

# Row 24
$ws.Range('D24').Value = '84.33'
$ws.Range('E24').Value = '  +3.03%  '

# Row 25
$ws.Range('E25').Value = '  +0.53%  '

# Row 26
$ws.Range('D26').Value = '12.42'
$ws.Range('E26').Value = '  +2.20%  '

# Row 27
$ws.Range('E27').Value = '  +3.88%  '

# Row 28
$ws.Range('E28').Value = '  +3.65%  '

# Row 29
$ws.Range('D29').Value = '3.955.18'
$ws.Range('E29').Value = '  +1.06%  '

# Row 30
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.11%  '

# Row 31
$ws.Range('D31').Value = '3.10'
$ws.Range('E31').Value = '  +15.55%  '

# Row 32
$ws.Range('E32').Value = '  +6.32%  '

# Row 33
$ws.Range('E33').Value = '  +1.84%  '

# Row 34
$ws.Range('D34').Value = '0.186'
$ws.Range('E34').Value = '  +7.00%  '

# Row 35
$ws.Range('D35').Value = '29.49'
$ws.Range('E35').Value = '  +3.68%  '

# Row 36
$ws.Range('E36').Value = '  +5.07%  '

# Row 37
$ws.Range('E37').Value = '  -0.04%  '

# Row 38
$ws.Range('E38').Value = '  +2.76%  '

# Row 39
$ws.Range('E39').Value = '  +7.15%  '

# Row 40
$ws.Range('E40').Value = '  +4.77%  '

# Row 41
$ws.Range('D41').Value = '2.22'
$ws.Range('E41').Value = '  +12.61%  '

# Row 42
$ws.Range('B42').Value = 'FLOKI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D42').Value = '0.000329'
$ws.Range('E42').Value = '  +23.93%  '

# Row 43
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value = '0.975'
$ws.Range('E43').Value = '  +2.12%  '

# Row 44
$ws.Range('E44').Value = '  -0.07%  '

# Row 46
$ws.Range('D46').Value = '162.36'
$ws.Range('E46').Value = '  +4.60%  '

# Row 47
$ws.Range('D47').Value = '45.03'
$ws.Range('E47').Value = '  +0.11%  '

# Row 48
$ws.Range('D48').Value = '48.87'
$ws.Range('E48').Value = '  +3.24%  '

# Row 49
$ws.Range('E49').Value = '  +2.61%  '

# Row 50
$ws.Range('E50').Value = '  -1.14%  '

# Row 51
$ws.Range('E51').Value = '  +2.89%  '
